$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "variable" header in C1 is now "key"
$ws.Range("C1").Value = "key"

# Update the selection to reflect what the author had selected when saving
$ws.Range("C2").Select()
